$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("个人持仓")
$ws1.Range("C2").Value = 50.33
$ws1.Range("C3").Value = 33.57
$ws1.Range("C4").Value = 48.58
$ws1.Range("C5").Value = 53.12
$ws1.Range("C6").Value = 30.92
$ws1.Range("C7").Value = 40.05
$ws1.Range("C8").Value = 26.29
$ws1.Range("C9").Value = 25.87
$ws1.Range("C10").Value = 25.04
$ws1.Range("C11").Value = 133.8
$ws1.Range("C12").Value = 159.92
$ws1.Range("C13").Value = 227
$ws1.Range("C14").Value = 0.801
$ws1.Range("C15").Value = 10.42
$ws1.Range("C16").Value = 26.6
$ws1.Range("C17").Value = 26.71
$ws1.Range("C18").Value = 16.05
$ws1.Range("C19").Value = 34.24
$ws1.Range("C20").Value = 36.71
$ws1.Range("C21").Value = 27.42
$ws1.Range("C22").Value = 109.3
$ws1.Range("C23").Value = 3.817
$ws1.Range("C24").Value = 4.561

$ws2 = $wb.Worksheets.Item("家庭持仓")
$ws2.Range("C2").Value = 50.33
$ws2.Range("C3").Value = 33.57
$ws2.Range("C4").Value = 48.58
$ws2.Range("C5").Value = 53.12
$ws2.Range("C6").Value = 30.92
$ws2.Range("C7").Value = 40.05
$ws2.Range("C8").Value = 26.29
$ws2.Range("C9").Value = 25.87
$ws2.Range("C10").Value = 25.04
$ws2.Range("C11").Value = 133.8
$ws2.Range("C12").Value = 227
$ws2.Range("C13").Value = 159.92
$ws2.Range("C14").Value = 246.31
$ws2.Range("C15").Value = 0.801
$ws2.Range("C16").Value = 39.3
$ws2.Range("C17").Value = 10.42
$ws2.Range("C18").Value = 26.6
$ws2.Range("C19").Value = 26.71
$ws2.Range("C20").Value = 16.05
$ws2.Range("C21").Value = 34.24
$ws2.Range("C22").Value = 36.71
$ws2.Range("C23").Value = 27.42
$ws2.Range("C24").Value = 109.3
$ws2.Range("C25").Value = 3.817
$ws2.Range("C26").Value = 1.1
$ws2.Range("C27").Value = 4.561

